$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.985.14"
$ws.Range("E2").Value = "  +2.62%  "

$ws.Range("D3").Value = "3.086.03"
$ws.Range("E3").Value = "  +4.67%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.12"
$ws.Range("E5").Value = "  +2.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.49"
$ws.Range("E6").Value = "  +5.63%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.079.18"
$ws.Range("E8").Value = "  +4.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  +1.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.61"
$ws.Range("E10").Value = "  -1.06%  "

$ws.Range("E11").Value = "  +2.55%  "

$ws.Range("E12").Value = "  +5.29%  "

$ws.Range("E13").Value = "  +1.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.50"
$ws.Range("E14").Value = "  +6.29%  "

$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("D16").Value = "3.597.63"
$ws.Range("E16").Value = "  +4.62%  "

$ws.Range("D17").Value = "66.961.90"
$ws.Range("E17").Value = "  +2.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  +3.37%  "

$ws.Range("D19").Value = "3.087.10"
$ws.Range("E19").Value = "  +3.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.18"
$ws.Range("E20").Value = "  +9.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "465.75"
$ws.Range("E21").Value = "  +4.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.716"
$ws.Range("E22").Value = "  +4.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("E23").Value = "  +4.08%  "

$ws.Range("E24").Value = "  +1.33%  "

$ws.Range("E25").Value = "  +6.65%  "

$ws.Range("E26").Value = "  +6.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  +1.19%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.96"
$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("E30").Value = "  +0.57%  "

$ws.Range("E31").Value = "  +3.38%  "

$ws.Range("E32").Value = "  +0.78%  "

$ws.Range("E33").Value = "  +3.44%  "

$ws.Range("E34").Value = "  +3.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("E36").Value = "  +2.90%  "

$ws.Range("E37").Value = "  +2.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.12"
$ws.Range("E38").Value = "  +7.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "46.81"
$ws.Range("E39").Value = "  +5.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.318"
$ws.Range("E40").Value = "  +6.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.17"
$ws.Range("E41").Value = "  +2.39%  "

$ws.Range("E42").Value = "  +1.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.68"
$ws.Range("E43").Value = "  +2.51%  "

$ws.Range("E44").Value = "  -0.73%  "

$ws.Range("E45").Value = "  +2.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "383.84"
$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("D47").Value = "2.763.16"
$ws.Range("E47").Value = "  +2.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.66"
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.64"
$ws.Range("E50").Value = "  +6.33%  "

$ws.Range("E51").Value = "  +2.80%  "
